$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Fund holdings detail for the new "2022-Q3" sheet.
# Columns: index, 基金代码, 基金名称, 基金规模, 股票总仓位, 仓位占比,
#          持有市值(亿元), 仓位排名
# ---------------------------------------------------------------------------
$fundData = @(
    @("0","002363","华安安康灵活配置混合A","107.76","21.95","1.03","1.1099","5"),
    @("1","009312","新疆前海联合价值优选混合A","8.45","92.63","6.57","0.5552","4"),
    @("2","004693","新疆前海联合泳隽灵活配置混合A","8.43","92.61","6.33","0.5336","4"),
    @("3","002364","华安安康灵活配置混合C","27.33","21.95","1.03","0.2815","5"),
    @("4","002144","华安新优选灵活配置混合C","19.79","22.13","1.17","0.2315","6"),
    @("5","001312","华安新优选灵活配置混合A","9.83","22.13","1.17","0.1150","6"),
    @("6","009313","新疆前海联合价值优选混合C","1.48","92.63","6.57","0.0972","4"),
    @("7","011390","华安添祥6个月持有期混合A","6.67","33.77","1.19","0.0794","8"),
    @("8","515400","富国中证大数据产业ETF","2.11","99.08","3.58","0.0755","10"),
    @("9","001028","华安物联网主题股票A","2.96","94.04","2.41","0.0713","5"),
    @("10","012321","东财中证云计算指数增强A","1.31","93.47","5.31","0.0696","5"),
    @("11","516510","易方达中证云计算与大数据主题ETF","1.70","98.74","3.36","0.0571","10"),
    @("12","161628","融通中证云计算与大数据主题指数（LOF）A","1.52","94.31","3.35","0.0509","10"),
    @("13","560660","新华中证云计算50ETF","1.20","97.03","3.86","0.0463","9"),
    @("14","970015","申万宏源红利成长灵活配置混合","0.86","80.01","3.56","0.0306","3"),
    @("15","012322","东财中证云计算指数增强C","0.51","93.47","5.31","0.0271","5"),
    @("16","516630","华夏中证云计算与大数据主题ETF","0.64","99.46","3.38","0.0216","10"),
    @("17","159739","鹏华中证云计算与大数据主题ETF","0.62","97.89","3.33","0.0206","10"),
    @("18","014130","融通中证云计算与大数据主题指数（LOF）C","0.58","94.31","3.35","0.0194","10"),
    @("19","516000","华夏中证大数据产业ETF","0.50","97.60","3.52","0.0176","10"),
    @("20","168701","合煦智远国证香蜜湖金融科技指数（LOF）A","0.47","92.07","2.63","0.0124","7"),
    @("21","159890","招商中证云计算与大数据主题ETF","0.36","98.52","3.34","0.0120","10"),
    @("22","005444","光大保德信多策略精选18个月定期开放灵活配置混合","0.84","29.28","1.41","0.0118","8"),
    @("23","000549","华安大国新经济股票A","0.72","93.83","1.46","0.0105","10"),
    @("24","516700","华宝中证大数据产业ETF","0.23","95.99","3.46","0.0080","10"),
    @("25","009409","华安添福18个月持有期混合A","0.72","21.56","1.10","0.0079","6"),
    @("26","090011","大成核心双动力混合","0.24","92.56","3.13","0.0075","8"),
    @("27","010307","西藏东财信息产业精选混合A","0.17","84.91","3.26","0.0055","10"),
    @("28","013132","创金合信文娱媒体股票型发起式证券投资基金A","0.06","91.90","6.18","0.0037","2"),
    @("29","168702","合煦智远国证香蜜湖金融科技指数（LOF）C","0.13","92.07","2.63","0.0034","7"),
    @("30","013133","创金合信文娱媒体股票型发起式证券投资基金C","0.05","91.90","6.18","0.0031","2"),
    @("31","010308","西藏东财信息产业精选混合C","0.08","84.91","3.26","0.0026","10"),
    @("32","009410","华安添福18个月持有期混合C","0.09","21.56","1.10","0.0010","6"),
    @("33","016291","华安大国新经济股票C","0.00","93.83","1.46","0.00","10"),
    @("34","007042","新疆前海联合泳隽灵活配置混合C","0.00","92.61","6.33","0.00","4"),
    @("35","016292","华安物联网主题股票C","0.00","94.04","2.41","0.00","5"),
    @("36","016181","华安添祥6个月持有期混合C","0.00","33.77","1.19","0.00","8")
)

# ---------------------------------------------------------------------------
# 1. Insert a new "2022-Q3" worksheet right after "总计" (before "2022-Q2").
# ---------------------------------------------------------------------------
$totalWs = $wb.Worksheets.Item(1)
$q3 = $wb.Worksheets.Add($null, $totalWs)
$q3.Name = "2022-Q3"

# ---------------------------------------------------------------------------
# 2. Update the "总计" summary sheet: insert a new row 2 for 2022-Q3 data and
#    renumber the index column for the rows that shifted down.
# ---------------------------------------------------------------------------
$totalWs.Rows.Item(2).Insert()

# Re-apply the formatting that the freshly inserted row lost by copying it
# from the row directly below (which still carries the original styling).
$totalWs.Range("A3:D3").Copy()
$totalWs.Range("A2:D2").PasteSpecial(-4122)

$totalWs.Range("A2").Value = 0
$totalWs.Range("B2").Value = "2022-Q3"
$totalWs.Range("C2").Value = 37
$totalWs.Range("D2").Value = 3.6

for ($i = 3; $i -le 9; $i++) {
    $totalWs.Cells.Item($i, 1).Value = $i - 2
}

# ---------------------------------------------------------------------------
# 3. Populate the new "2022-Q3" sheet with the header row and fund holdings.
# ---------------------------------------------------------------------------

# Seed the bold/border/center-top formatting (style used by every other
# quarter sheet header + index column) by copying it from "总计"!A2.
$totalWs.Range("A2").Copy()
$q3.Range("B1:H1").PasteSpecial(-4122)
$q3.Range("A2:A38").PasteSpecial(-4122)

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($col = 2; $col -le 8; $col++) {
    $q3.Cells.Item(1, $col).Value = $headers[$col - 2]
}

# Force 基金代码/基金规模/股票总仓位/仓位占比/持有市值 to be stored as text so
# leading zeros and the source decimal formatting are preserved verbatim.
$q3.Range("B2:B38").NumberFormat = "@"
$q3.Range("D2:G34").NumberFormat = "@"

foreach ($row in $fundData) {
    $r = [int]$row[0] + 2

    $q3.Cells.Item($r, 1).Value = [int]$row[0]
    $q3.Cells.Item($r, 2).Value = $row[1]
    $q3.Cells.Item($r, 3).Value = $row[2]
    $q3.Cells.Item($r, 4).Value = $row[3]
    $q3.Cells.Item($r, 5).Value = $row[4]
    $q3.Cells.Item($r, 6).Value = $row[5]

    if ($r -le 34) {
        $q3.Cells.Item($r, 7).Value = $row[6]
    } else {
        $q3.Cells.Item($r, 7).Value = 0
    }

    $q3.Cells.Item($r, 8).Value = [int]$row[7]
}

Write-Host "done"
